$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (Volume number, date range) ---
$a8 = $ws.Range("A8")
$a8.Characters(21,1).Text = "10"

$c9 = $ws.Range("C9")
$c9.Characters(47,8).Text = "3/9/2025"
$c9.Characters(27,9).Text = "3/3/2025"

# --- Cells that only need a value change (style unchanged) ---
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -53.846153846153
$ws.Range("J16").Value = 27
$ws.Range("K16").Value = -55.555555555555
$ws.Range("L16").Value = -69.230769230769
$ws.Range("M16").Value = -61.290322580645
$ws.Range("N16").Value = -92.356687898089
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 14
$ws.Range("K17").Value = -26.315789473684
$ws.Range("L17").Value = -51.724137931034
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = -54.838709677419
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -40
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -45
$ws.Range("I18").Value = 27
$ws.Range("J18").Value = 46
$ws.Range("K18").Value = -41.304347826087
$ws.Range("L18").Value = -50
$ws.Range("M18").Value = -41.304347826087
$ws.Range("N18").Value = -79.389312977099
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 70
$ws.Range("G19").Value = 72
$ws.Range("H19").Value = -2.777777777777
$ws.Range("I19").Value = 161
$ws.Range("J19").Value = 187
$ws.Range("K19").Value = -13.903743315508
$ws.Range("L19").Value = -33.195020746888
$ws.Range("M19").Value = -7.471264367816
$ws.Range("N19").Value = -60.731707317073
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -60
$ws.Range("I20").Value = 3
$ws.Range("J20").Value = 8
$ws.Range("K20").Value = -62.5
$ws.Range("L20").Value = -40
$ws.Range("M20").Value = -40
$ws.Range("N20").Value = -98.064516129032
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -40.625
$ws.Range("G21").Value = 120
$ws.Range("H21").Value = -20.833333333333
$ws.Range("I21").Value = 218
$ws.Range("J21").Value = 288
$ws.Range("K21").Value = -24.305555555555
$ws.Range("L21").Value = -40.921409214092
$ws.Range("M21").Value = -19.557195571955
$ws.Range("N21").Value = -75.395033860045
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -16.666666666666
$ws.Range("J22").Value = 11
$ws.Range("K22").Value = -9.090909090909
$ws.Range("L22").Value = -23.076923076923
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -20.689655172413
$ws.Range("F24").Value = 113
$ws.Range("G24").Value = 128
$ws.Range("H24").Value = -11.71875
$ws.Range("I24").Value = 266
$ws.Range("J24").Value = 298
$ws.Range("K24").Value = -10.738255033557
$ws.Range("L24").Value = -17.133956386292
$ws.Range("M24").Value = 0.377358490566
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = -26.923076923076
$ws.Range("G25").Value = 111
$ws.Range("H25").Value = -18.918918918918
$ws.Range("I25").Value = 198
$ws.Range("J25").Value = 248
$ws.Range("K25").Value = -20.161290322580
$ws.Range("L25").Value = -16.806722689075
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -10
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = 3.571428571428
$ws.Range("I26").Value = 62
$ws.Range("J26").Value = 66
$ws.Range("K26").Value = -6.060606060606
$ws.Range("L26").Value = -11.428571428571
$ws.Range("M26").Value = 82.352941176470
$ws.Range("G27").Value = 2
$ws.Range("J27").Value = 2
$ws.Range("K27").Value = 0
$ws.Range("C28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 10
$ws.Range("H28").Value = -10
$ws.Range("I28").Value = 17
$ws.Range("J28").Value = 14
$ws.Range("K28").Value = 21.428571428571
$ws.Range("L28").Value = 0
$ws.Range("I31").Value = 3
$ws.Range("K31").Value = 200

# --- Cells that need both a style (number-format) change and a value change ---
# Donor cells (never themselves edited) supply number formats via PasteSpecial(xlPasteFormats=-4122):
#   $ws.Range("C14") -> style 13 (text/no format)
#   $ws.Range("D16") -> style 14 (integer #,##0)
#   $ws.Range("L15") -> style 15 (percent #,##0.0)
$ws.Range("D16").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1

$ws.Range("L15").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100

$ws.Range("D16").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1

$ws.Range("L15").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = -100

$ws.Range("D16").Copy()
$ws.Range("J15").PasteSpecial(-4122)
$ws.Range("J15").Value = 1

$ws.Range("L15").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("K15").Value = 0

$ws.Range("C16").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("D17").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D17").PasteSpecial(-4122)

$ws.Range("E17").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("D16").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 1

$ws.Range("C22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D16").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("F31").Value = 1
